$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Select entire row 9 and insert a new blank row above it (shifts rows 9:200 down to 10:201)
$rng = $ws.Range("9:9")
$rng.Select()
$rng.Insert()

# Restore the explicit row heights that Excel's autofit recomputed for the
# wrap-text "category" rows affected by the shift.
$ws.Rows.Item(8).RowHeight = 68
$ws.Rows.Item(57).RowHeight = 46
$ws.Rows.Item(76).RowHeight = 41
$ws.Rows.Item(85).RowHeight = 68
$ws.Rows.Item(101).RowHeight = 86
$ws.Rows.Item(108).RowHeight = 71
